$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 983.65
$ws.Range("I28").Value = 858.125
$ws.Range("K28").Value = 858.125
$ws.Range("M28").Value = -373.125

$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H113").Value = 8953.764999999999
$ws.Range("J113").Value = 9288.666999999999
$ws.Range("L113").Value = 9288.666999999999
$ws.Range("N113").Value = -15796.667

$ws.Range("H138").Value = 2989.55
$ws.Range("I138").Value = 1640.05
$ws.Range("J138").Value = 3664.3
$ws.Range("K138").Value = 4920.15
$ws.Range("L138").Value = 10992.9
$ws.Range("M138").Value = 219.8500000000004
$ws.Range("N138").Value = -21272.9

$ws.Range("H141").Value = 4869.1924
$ws.Range("I141").Value = 3825.4092
$ws.Range("J141").Value = 10610
$ws.Range("K141").Value = 11476.2276
$ws.Range("L141").Value = 31830
$ws.Range("M141").Value = -6296.2276
$ws.Range("N141").Value = -42190

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 25019.715
$ws.Range("I2").Value = 47381.285
$ws.Range("K2").Value = 47381.285
$ws.Range("M2").Value = -47268.285

$ws.Range("H45").Value = 5234.1113
$ws.Range("I45").Value = 1902
$ws.Range("J45").Value = 7899.8
$ws.Range("K45").Value = 1902
$ws.Range("L45").Value = 7899.8
$ws.Range("M45").Value = -1525
$ws.Range("N45").Value = -8653.799999999999

$ws.Range("H61").Value = 1608.125
$ws.Range("I61").Value = 1460.6522
$ws.Range("K61").Value = 1460.6522
$ws.Range("M61").Value = -1248.6522

$ws.Range("H102").Value = 2996.55
$ws.Range("I102").Value = 3101.6316
$ws.Range("K102").Value = 3101.6316
$ws.Range("M102").Value = -1479.6316

$ws.Range("H116").Value = 25019.715
$ws.Range("I116").Value = 47381.285
$ws.Range("K116").Value = 47381.285
$ws.Range("M116").Value = -45087.285

$ws.Range("H136").Value = 1608.125
$ws.Range("I136").Value = 1460.6522
$ws.Range("K136").Value = 4381.9566
$ws.Range("M136").Value = -1831.9566

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 25019.715
$ws.Range("I3").Value = 47381.285
$ws.Range("K3").Value = 47381.285
$ws.Range("M3").Value = -47267.285

$ws.Range("H107").Value = 4307.7095
$ws.Range("I107").Value = 4005.68
$ws.Range("J107").Value = 5566.1665
$ws.Range("K107").Value = 4005.68
$ws.Range("L107").Value = 5566.1665
$ws.Range("M107").Value = -2085.68
$ws.Range("N107").Value = -9406.166499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 148.79167
$ws.Range("I7").Value = 28.071428
$ws.Range("J7").Value = 317.8
$ws.Range("K7").Value = 28.071428
$ws.Range("L7").Value = 317.8
$ws.Range("M7").Value = 84.928572
$ws.Range("N7").Value = -543.8

$ws.Range("H16").Value = 2856.818
$ws.Range("I16").Value = 2658.5
$ws.Range("K16").Value = 2658.5
$ws.Range("M16").Value = -2371.5

$ws.Range("H41").Value = 3639.0557
$ws.Range("I41").Value = 3718.8125
$ws.Range("J41").Value = 3001
$ws.Range("K41").Value = 3718.8125
$ws.Range("L41").Value = 3001
$ws.Range("M41").Value = -3290.8125
$ws.Range("N41").Value = -3857

$ws.Range("H51").Value = 10333
$ws.Range("I51").Value = 10333
$ws.Range("K51").Value = 10333
$ws.Range("M51").Value = -9597

$ws.Range("H59").Value = 999999
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws.Range("H60").Value = 632.2174
$ws.Range("I60").Value = 276.6
$ws.Range("J60").Value = 3003
$ws.Range("K60").Value = 276.6
$ws.Range("L60").Value = 3003
$ws.Range("M60").Value = 234.4
$ws.Range("N60").Value = -4025

$ws.Range("H61").Value = 10333
$ws.Range("I61").Value = 10333
$ws.Range("K61").Value = 10333
$ws.Range("M61").Value = -9985

$ws.Range("H113").Value = 2856.818
$ws.Range("I113").Value = 2658.5
$ws.Range("K113").Value = 2658.5
$ws.Range("M113").Value = -488.5

$ws.Range("H118").Value = 27500
$ws.Range("J118").Value = 27500
$ws.Range("L118").Value = 27500
$ws.Range("N118").Value = -30814

$ws.Range("H122").Value = 3834.9707
$ws.Range("I122").Value = 2974.5
$ws.Range("J122").Value = 5064.2144
$ws.Range("K122").Value = 8923.5
$ws.Range("L122").Value = 15192.6432
$ws.Range("M122").Value = -6473.5
$ws.Range("N122").Value = -20092.6432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 3
$ws.Range("I12").Value = 1
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 170

$ws.Range("H113").Value = 2504
$ws.Range("J113").Value = 1889.3334
$ws.Range("L113").Value = 5668.0002
$ws.Range("N113").Value = -10008.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3350.9
$ws.Range("I102").Value = 2556.5557
$ws.Range("J102").Value = 10500
$ws.Range("K102").Value = 2556.5557
$ws.Range("L102").Value = 10500
$ws.Range("M102").Value = -934.5556999999999
$ws.Range("N102").Value = -13744

$ws.Range("H107").Value = 777.1111
$ws.Range("I107").Value = 580.75
$ws.Range("K107").Value = 580.75
$ws.Range("M107").Value = 1339.25

$ws.Range("H126").Value = 3536.5
$ws.Range("J126").Value = 2849.5
$ws.Range("L126").Value = 8548.5
$ws.Range("N126").Value = -13488.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2498.3333
$ws.Range("I61").Value = 2435.75
$ws.Range("K61").Value = 2435.75
$ws.Range("M61").Value = -2233.75

$ws.Range("H113").Value = 2498.3333
$ws.Range("I113").Value = 2435.75
$ws.Range("K113").Value = 2435.75
$ws.Range("M113").Value = -265.75

$ws.Range("H122").Value = 7151.6
$ws.Range("I122").Value = 6752.6665
$ws.Range("K122").Value = 20257.9995
$ws.Range("M122").Value = -17807.9995

$ws.Range("H132").Value = 1867.1428
$ws.Range("I132").Value = 1627.9412
$ws.Range("K132").Value = 4883.8236
$ws.Range("M132").Value = -2353.8236

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 9375
$ws.Range("J20").Value = 4000
$ws.Range("L20").Value = 4000
$ws.Range("N20").Value = -4480

$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()

$ws.Range("H29").Value = 7836.6665
$ws.Range("I29").Value = 9255
$ws.Range("K29").Value = 9255
$ws.Range("M29").Value = -8965

$ws.Range("H113").Value = 1349.6
$ws.Range("I113").Value = 974.5
$ws.Range("K113").Value = 2923.5
$ws.Range("M113").Value = -753.5

$ws.Range("H122").Value = 2243.4285
$ws.Range("I122").Value = 2243.4285
$ws.Range("K122").Value = 6730.2855
$ws.Range("M122").Value = -4280.2855

$ws.Range("H125").Value = 69999.82000000001
$ws.Range("J125").Value = 69999.82000000001
$ws.Range("L125").Value = 69999.82000000001
$ws.Range("N125").Value = -79839.82000000001

$ws.Range("H126").Value = 5428.75
$ws.Range("I126").Value = 6124.5
$ws.Range("J126").Value = 1950
$ws.Range("K126").Value = 18373.5
$ws.Range("L126").Value = 5850
$ws.Range("M126").Value = -15903.5
$ws.Range("N126").Value = -10790

$ws.Range("H128").Value = 51266.2
$ws.Range("J128").Value = 51266.2
$ws.Range("L128").Value = 51266.2
$ws.Range("N128").Value = -61226.2

$ws.Range("H129").Value = 64427.715
$ws.Range("J129").Value = 64427.715
$ws.Range("L129").Value = 64427.715
$ws.Range("N129").Value = -74427.715
